# Auto-generated Excel COM-interop script
# Updates market-board derived value columns (H-N) on several worksheets
# to match a scheduled data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
    # Row 43
    $ws.Range("H43").Value = 4036.25
    $ws.Range("I43").Value = 4036.25
    $ws.Range("K43").Value = 4036.25
    $ws.Range("M43").Value = -3967.25
    # Row 92
    $ws.Range("H92").Value = 2696.6667
    $ws.Range("I92").Value = 1656.7273
    $ws.Range("J92").Value = 5556.5
    $ws.Range("K92").Value = 1656.7273
    $ws.Range("L92").Value = 5556.5
    $ws.Range("M92").Value = -408.7273
    $ws.Range("N92").Value = -8052.5
    # Row 98
    $ws.Range("H98").Value = 2842651.5
    $ws.Range("I98").Value = 3126366.5
    $ws.Range("J98").Value = 5502.75
    $ws.Range("K98").Value = 3126366.5
    $ws.Range("L98").Value = 5502.75
    $ws.Range("M98").Value = -3124868.5
    $ws.Range("N98").Value = -8498.75
    # Row 107
    $ws.Range("H107").Value = 1741.3572
    $ws.Range("I107").Value = 438.5
    $ws.Range("J107").Value = 4998.5
    $ws.Range("K107").Value = 438.5
    $ws.Range("L107").Value = 4998.5
    $ws.Range("M107").Value = 1481.5
    $ws.Range("N107").Value = -8838.5
    # Row 116
    $ws.Range("H116").Value = 12789.19
    $ws.Range("I116").Value = 6447.222
    $ws.Range("J116").Value = 17545.666
    $ws.Range("K116").Value = 6447.222
    $ws.Range("L116").Value = 17545.666
    $ws.Range("M116").Value = -3005.222
    $ws.Range("N116").Value = -24429.666
    # Row 122
    $ws.Range("H122").Value = 2842651.5
    $ws.Range("I122").Value = 3126366.5
    $ws.Range("J122").Value = 5502.75
    $ws.Range("K122").Value = 9379099.5
    $ws.Range("L122").Value = 16508.25
    $ws.Range("M122").Value = -9376649.5
    $ws.Range("N122").Value = -21408.25
    # Row 137
    $ws.Range("H137").Value = 33336156
    $ws.Range("J137").Value = 3571.2856
    $ws.Range("L137").Value = 10713.8568
    $ws.Range("N137").Value = -15813.8568

$ws = $wb.Worksheets.Item("ARM")
    # Row 61
    $ws.Range("H61").Value = 13658709
    $ws.Range("I61").Value = 15913594
    $ws.Range("K61").Value = 15913594
    $ws.Range("M61").Value = -15913382
    # Row 74
    $ws.Range("H74").Value = 2884.1155
    $ws.Range("I74").Value = 3016.524
    $ws.Range("J74").Value = 2328
    $ws.Range("K74").Value = 3016.524
    $ws.Range("L74").Value = 2328
    $ws.Range("M74").Value = -2142.524
    $ws.Range("N74").Value = -4076
    # Row 77
    $ws.Range("H77").Value = 2884.1155
    $ws.Range("I77").Value = 3016.524
    $ws.Range("J77").Value = 2328
    $ws.Range("K77").Value = 15082.62
    $ws.Range("L77").Value = 11640
    $ws.Range("M77").Value = -10714.62
    $ws.Range("N77").Value = -20376
    # Row 97
    $ws.Range("H97").Value = 1032.7273
    $ws.Range("I97").Value = 1157.7084
    $ws.Range("J97").Value = 699.44446
    $ws.Range("K97").Value = 1157.7084
    $ws.Range("L97").Value = 699.44446
    $ws.Range("M97").Value = -661.7084
    $ws.Range("N97").Value = -1691.44446
    # Row 102
    $ws.Range("H102").Value = 1344.7646
    $ws.Range("I102").Value = 1303.8125
    $ws.Range("K102").Value = 1303.8125
    $ws.Range("M102").Value = 318.1875
    # Row 110
    $ws.Range("H110").Value = 5661.3447
    $ws.Range("I110").Value = 5653.4165
    $ws.Range("J110").Value = 5699.4
    $ws.Range("K110").Value = 5653.4165
    $ws.Range("L110").Value = 5699.4
    $ws.Range("M110").Value = -3608.4165
    $ws.Range("N110").Value = -9789.4
    # Row 136
    $ws.Range("H136").Value = 13658709
    $ws.Range("I136").Value = 15913594
    $ws.Range("K136").Value = 47740782
    $ws.Range("M136").Value = -47738232

$ws = $wb.Worksheets.Item("BSM")
    # Row 105
    $ws.Range("H105").Value = 849620.5600000001
    $ws.Range("I105").Value = 1040657.7
    $ws.Range("K105").Value = 1040657.7
    $ws.Range("M105").Value = -1038910.7

$ws = $wb.Worksheets.Item("CRP")
    # Row 7
    $ws.Range("H7").Value = 13.071428
    $ws.Range("I7").Value = 11.75
    $ws.Range("J7").Value = 21
    $ws.Range("K7").Value = 11.75
    $ws.Range("L7").Value = 21
    $ws.Range("M7").Value = 101.25
    $ws.Range("N7").Value = -247
    # Row 31
    $ws.Range("H31").Value = 19233776
    $ws.Range("I31").Value = 26317980
    $ws.Range("K31").Value = 26317980
    $ws.Range("M31").Value = -26317685
    # Row 34
    $ws.Range("H34").Value = 19233776
    $ws.Range("I34").Value = 26317980
    $ws.Range("K34").Value = 26317980
    $ws.Range("M34").Value = -26317778
    # Row 134
    $ws.Range("H134").Value = 1768.5588
    $ws.Range("I134").Value = 1998.1904
    $ws.Range("J134").Value = 1397.6154
    $ws.Range("K134").Value = 5994.5712
    $ws.Range("L134").Value = 4192.8462
    $ws.Range("M134").Value = -3459.5712
    $ws.Range("N134").Value = -9262.8462

$ws = $wb.Worksheets.Item("CUL")
    # Row 2
    $ws.Range("H2").Value = 203.54546
    $ws.Range("J2").Value = 187.25
    $ws.Range("L2").Value = 1123.5
    $ws.Range("N2").Value = -1349.5
    # Row 94
    $ws.Range("H94").Value = 10613.667
    $ws.Range("J94").Value = 18866.2
    $ws.Range("L94").Value = 56598.60000000001
    $ws.Range("N94").Value = -57950.60000000001
    # Row 98
    $ws.Range("H98").Value = 726.1539
    $ws.Range("J98").Value = 671.55554
    $ws.Range("L98").Value = 2014.66662
    $ws.Range("N98").Value = -5010.66662
    # Row 113
    $ws.Range("H113").Value = 1304
    $ws.Range("J113").Value = 1210
    $ws.Range("L113").Value = 3630
    $ws.Range("N113").Value = -7970
    # Row 117
    $ws.Range("H117").Value = 6066.273
    $ws.Range("I117").Value = 574.75
    $ws.Range("J117").Value = 9204.286
    $ws.Range("K117").Value = 1724.25
    $ws.Range("L117").Value = 27612.858
    $ws.Range("M117").Value = 1717.75
    $ws.Range("N117").Value = -34496.858
    # Row 132
    $ws.Range("H132").Value = 2617.6
    $ws.Range("I132").Value = 2129.1428
    $ws.Range("J132").Value = 3757.3333
    $ws.Range("K132").Value = 19162.2852
    $ws.Range("L132").Value = 33815.9997
    $ws.Range("M132").Value = -16632.2852
    $ws.Range("N132").Value = -38875.9997

$ws = $wb.Worksheets.Item("GSM")
    # Row 122
    $ws.Range("H122").Value = 5474.8335
    $ws.Range("I122").Value = 3579
    $ws.Range("J122").Value = 6829
    $ws.Range("K122").Value = 10737
    $ws.Range("L122").Value = 20487
    $ws.Range("M122").Value = -8287
    $ws.Range("N122").Value = -25387

$ws = $wb.Worksheets.Item("LTW")
    # Row 13
    $ws.Range("H13").Value = 0
    $ws.Range("I13").Value = 0
    $ws.Range("J13").Value = 0
    $ws.Range("K13").Value = 0
    $ws.Range("L13").Value = 0
    $ws.Range("M13:N13").ClearContents()
    # Row 40
    $ws.Range("H40").Value = 6676.1113
    $ws.Range("I40").Value = 6635.625
    $ws.Range("K40").Value = 6635.625
    $ws.Range("M40").Value = -6499.625
    # Row 100
    $ws.Range("H100").Value = 19253510
    $ws.Range("I100").Value = 3603
    $ws.Range("J100").Value = 25028482
    $ws.Range("K100").Value = 3603
    $ws.Range("L100").Value = 25028482
    $ws.Range("M100").Value = -3062
    $ws.Range("N100").Value = -25029564
    # Row 122
    $ws.Range("H122").Value = 3669
    $ws.Range("I122").Value = 3174.1562
    $ws.Range("J122").Value = 5648.375
    $ws.Range("K122").Value = 9522.4686
    $ws.Range("L122").Value = 16945.125
    $ws.Range("M122").Value = -7072.4686
    $ws.Range("N122").Value = -21845.125
    # Row 132
    $ws.Range("H132").Value = 3360.0571
    $ws.Range("I132").Value = 2282.6875
    $ws.Range("K132").Value = 6848.0625
    $ws.Range("M132").Value = -4318.0625
    # Row 136
    $ws.Range("H136").Value = 6468.5454
    $ws.Range("I136").Value = 6505.4614
    $ws.Range("J136").Value = 6415.222
    $ws.Range("K136").Value = 19516.3842
    $ws.Range("L136").Value = 19245.666
    $ws.Range("M136").Value = -16966.3842
    $ws.Range("N136").Value = -24345.666

$ws = $wb.Worksheets.Item("WVR")
    # Row 126
    $ws.Range("H126").Value = 5129.136
    $ws.Range("I126").Value = 5869.727
    $ws.Range("J126").Value = 4388.5454
    $ws.Range("K126").Value = 17609.181
    $ws.Range("L126").Value = 13165.6362
    $ws.Range("M126").Value = -15139.181
    $ws.Range("N126").Value = -18105.6362
    # Row 132
    $ws.Range("H132").Value = 264569.8
    $ws.Range("I132").Value = 1007.37036
    $ws.Range("J132").Value = 911495.8
    $ws.Range("K132").Value = 3022.11108
    $ws.Range("L132").Value = 2734487.4
    $ws.Range("M132").Value = -492.1110800000001
    $ws.Range("N132").Value = -2739547.4
    # Row 136
    $ws.Range("H136").Value = 313466.28
    $ws.Range("I136").Value = 981.4583
    $ws.Range("K136").Value = 2944.3749
    $ws.Range("M136").Value = -394.3748999999998
